$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 (GNN-MT), shifting existing
# GNN-MT/RF/PN rows down by one (to rows 3/4/5).
$ws.Rows.Item(2).Insert()

# The row insert leaves B2 with inherited formatting from the header row;
# this new cell should be plain (no style), matching the other value cells.
$ws.Range("B2").ClearFormats()

# New row 2: GNN-MT-O. Copy A3's formatting (bold/border/centered) onto A2,
# then set its text and value.
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("A2").Value = "16_train (GNN-MT-O) val delta-auprc"
$ws.Range("B2").Value = 3.644117647058823

# Update the values of the shifted rows (text/formatting stay the same)
$ws.Range("B3").Value = 3.252941176470588
$ws.Range("B4").Value = 2.920588235294117
$ws.Range("B5").Value = 2.614705882352941

# New row 6: PN-O, appended after PN. Copy A5's formatting onto A6.
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("A6").Value = "16_train (PN-O) val delta-auprc"
$ws.Range("B6").Value = 2.56764705882353
